# Apply the edit described in the diff:
#  1. Insert a new worksheet "Player Info" before the existing "ODI Batting"
#     sheet, with ID / NAME / BATTING_HAND / BOWL_STYLE columns.
#  2. Rename the "MATCH_CARD_LINK" column to "MATCH_CODE" on both the
#     "ODI Batting" and "ODI Bowling" sheets, and replace the full
#     howstat.com URLs in that column with just the bare match code number.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Add the new "Player Info" sheet in front of "ODI Batting"
# ---------------------------------------------------------------------
# NOTE: worksheet handles returned by this COM runtime track a *position*
# rather than a specific sheet, so grab the "before" sheet, insert, rename
# and then re-fetch the other sheets (by name) fresh afterwards - any
# handle obtained before the insert may now point at a different sheet.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$infoSheet = $wb.Worksheets.Add($battingSheet)
$infoSheet.Name = "Player Info"

$header = $infoSheet.Range("A1:D1")
$header.Font.Bold = $true
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4160
$header.Borders.LineStyle = 1
$header.Borders.Weight = 2

$infoSheet.Range("A1").Value = "ID"
$infoSheet.Range("B1").Value = "NAME"
$infoSheet.Range("C1").Value = "BATTING_HAND"
$infoSheet.Range("D1").Value = "BOWL_STYLE"

# Player id looks numeric - force text storage to match source data style
$infoSheet.Range("A2").NumberFormat = "@"
$infoSheet.Range("A2").Value = "4566"
$infoSheet.Range("B2").Value = "Liam A Dawson"
$infoSheet.Range("C2").Value = "Right Handed"
$infoSheet.Range("D2").Value = "Left Arm Orthodox"

# ---------------------------------------------------------------------
# 2. Rename MATCH_CARD_LINK -> MATCH_CODE and replace urls with codes
# ---------------------------------------------------------------------

# Re-fetch sheets by name now that the sheet collection has been
# reordered by the insert above.
$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# "ODI Batting" sheet - the link lives in column D
$battingSheet.Range("D1").Value = "MATCH_CODE"

$battingCodes = @{
    2 = "3932"
    3 = "4209"
    4 = "4210"
    5 = "4660"
    6 = "4663"
    7 = "4666"
}
foreach ($row in $battingCodes.Keys) {
    $cell = $battingSheet.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $battingCodes[$row]
}

# "ODI Bowling" sheet - the link lives in column B
$bowlingSheet.Range("B1").Value = "MATCH_CODE"

$bowlingCodes = @{
    2 = "3932"
    3 = "4210"
    4 = "4660"
    5 = "4663"
    6 = "4666"
}
foreach ($row in $bowlingCodes.Keys) {
    $cell = $bowlingSheet.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $bowlingCodes[$row]
}
